# Commit: Sat, Apr 04, 2020  6:05:43 PM
#
# 1) Slide 16's table switches from the deck's custom "Table_0" style to the
#    built-in PowerPoint "Table Grid" style.
# 2) The deck's theme colours change from the custom "Integral" palette to
#    the standard Office palette (12 theme colours: dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 ---------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{A575934A-081C-438B-8B24-333417AFC5B9}")

# --- 2) Theme colours: Integral -> Office ------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme
# Order matches the OOXML <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = $hex -band 0xFF0000
    $r = [math]::Floor($r / 0x10000)
    $g = $hex -band 0x00FF00
    $g = [math]::Floor($g / 0x100)
    $b = $hex -band 0x0000FF
    # RGB COM values are packed 0xBBGGRR (i.e. R + G*256 + B*65536)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
